# OM_EM_Scenarios_v2.xlsx -- "update simulation runs + plots"
#
# 1. OM sheet (sheet1 / "OM"): tweak several scenario-parameter columns on
#    rows 2-5 (J,N,O,P,S,T,U,V,W,X,Y,Z).
# 2. EM_1Fl_RW sheet (sheet4 / "EM_1Fl_RW"): append 18 new scenario rows
#    (38-55) covering the new RW_1.75 and RW_Est scenarios.
# 3. View-state housekeeping: move the active tab from OM to EM_1Fl_RW, and
#    update each sheet's remembered selection.

$wb = $excel.ActiveWorkbook
$wsOM  = $wb.Worksheets.Item(1)   # OM
$wsRW  = $wb.Worksheets.Item(4)   # EM_1Fl_RW
$wsST  = $wb.Worksheets.Item(5)   # EM_Fast_Blk_SensTest

# ---------------------------------------------------------------------
# 1. OM sheet parameter edits (rows 2-5)
# ---------------------------------------------------------------------
$omEdits = @{
    2 = @{ J=0.85; N=0.2; O=15; P=15; S=5.5;  T=1.75; U=8.5; V=2.25; W=9.5;  X=1.25;               Y=12.5; Z=1.25 }
    3 = @{ J=0.85; N=0.2; O=15; P=15; S=5.5;  T=1.75; U=8.5; V=2.25; W=17;   X=8;                  Y=20;   Z=7    }
    4 = @{ J=0.85; N=0.2; O=15; P=15; S=5.5;  T=1.75; U=8.5; V=2.25; W=9.5;  X=1.1499999999999999; Y=12.5; Z=1.25 }
    5 = @{ J=0.85; N=0.2; O=15; P=15; S=5.5;  T=1.75; U=8.5; V=2.25; W=17;   X=8;                  Y=20;   Z=7    }
}

foreach ($r in $omEdits.Keys) {
    $cols = $omEdits[$r]
    foreach ($col in $cols.Keys) {
        $wsOM.Range("$col$r").Value = $cols[$col]
    }
}

# ---------------------------------------------------------------------
# 2. EM_1Fl_RW sheet: append rows 38-55
# ---------------------------------------------------------------------
$rwRows = @(
    @("Term_1Fl_Gam_RW_1.75",  1, "50,70", "gamma",        "RW", 2, "ln_fish_selpars_re", $false, 1.75),
    @("TrxE_1Fl_Gam_RW_1.75",  1, "30,50", "gamma",        "RW", 2, "ln_fish_selpars_re", $false, 1.75),
    @("Int_1Fl_Gam_RW_1.75",   1, "28,39", "gamma",        "RW", 2, "ln_fish_selpars_re", $false, 1.75),
    @("Term_1Fl_L_RW_1.75",    1, "50,70", "logistic",     "RW", 2, "ln_fish_selpars_re", $false, 1.75),
    @("TrxE_1Fl_L_RW_1.75",    1, "30,50", "logistic",     "RW", 2, "ln_fish_selpars_re", $false, 1.75),
    @("Int_1Fl_L_RW_1.75",     1, "28,39", "logistic",     "RW", 2, "ln_fish_selpars_re", $false, 1.75),
    @("Term_1Fl_ExpL_RW_1.75", 1, "50,70", "exp_logistic", "RW", 3, "ln_fish_selpars_re", $false, 1.75),
    @("TrxE_1Fl_ExpL_RW_1.75", 1, "30,50", "exp_logistic", "RW", 3, "ln_fish_selpars_re", $false, 1.75),
    @("Int_1Fl_ExpL_RW_1.75",  1, "28,39", "exp_logistic", "RW", 3, "ln_fish_selpars_re", $false, 1.75),
    @("Term_1Fl_ExpL_RW_Est",  1, "50,70", "exp_logistic", "RW", 3, "ln_fish_selpars_re", $false, "NA"),
    @("Term_1Fl_L_RW_Est",     1, "50,70", "logistic",     "RW", 2, "ln_fish_selpars_re", $false, "NA"),
    @("Term_1Fl_Gam_RW_Est",   1, "50,70", "gamma",        "RW", 2, "ln_fish_selpars_re", $false, "NA"),
    @("TrxE_1Fl_ExpL_RW_Est",  1, "30,50", "exp_logistic", "RW", 3, "ln_fish_selpars_re", $false, "NA"),
    @("TrxE_1Fl_L_RW_Est",     1, "30,50", "logistic",     "RW", 2, "ln_fish_selpars_re", $false, "NA"),
    @("TrxE_1Fl_Gam_RW_Est",   1, "30,50", "gamma",        "RW", 2, "ln_fish_selpars_re", $false, "NA"),
    @("Int_1Fl_ExpL_RW_Est",   1, "28,39", "exp_logistic", "RW", 3, "ln_fish_selpars_re", $false, "NA"),
    @("Int_1Fl_L_RW_Est",      1, "28,39", "logistic",     "RW", 2, "ln_fish_selpars_re", $false, "NA"),
    @("Int_1Fl_Gam_RW_Est",    1, "28,39", "gamma",        "RW", 2, "ln_fish_selpars_re", $false, "NA")
)

$startRow = 38
for ($i = 0; $i -lt $rwRows.Length; $i++) {
    $row = $rwRows[$i]
    $r = $startRow + $i
    for ($c = 0; $c -lt $row.Length; $c++) {
        $wsRW.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}

# ---------------------------------------------------------------------
# 3. View-state: selections + active tab
# ---------------------------------------------------------------------
# OM keeps the D6 selection no longer; new remembered selection is W6, and
# it is no longer the front-most tab.
$wsOM.Select()
$excel.ActiveWindow.ScrollColumn = 14
$excel.ActiveWindow.ScrollRow = 1
$wsOM.Range("W6").Select()

# EM_Fast_Blk_SensTest: selection collapses from A3:G16 to the single cell E21.
$wsST.Range("E21").Select()

# EM_1Fl_RW becomes the active tab, scrolled down to the newly-added rows,
# with the new selection on C55.
$wsRW.Select()
$excel.ActiveWindow.ScrollRow = 23
$excel.ActiveWindow.ScrollColumn = 1
$wsRW.Range("C55").Select()
